# Adding data observability code
# Applies:
#  1. Updated drift statistics on "Drift" sheet
#  2. New "skewness" column + refreshed stats on "Numerical" sheet
#  3. Refreshed row counts on "Categorical" sheet
#  4. New "Usage" sheet with table/index usage info

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Drift sheet - refresh mean/std drift numbers (columns B & C, rows 2-8)
# ---------------------------------------------------------------------
$drift = $wb.Worksheets.Item("Drift")

$driftValues = @(
    @(2, 5009.642359663751, 808.0612737739293),
    @(3, 263.4681858581362, 1971.477389070604),
    @(4, 137.170842047112,  791.8762270096495),
    @(5, 1755.276375161149, 3799.524769506221),
    @(6, 1056.006506297232, 1779.802687801041),
    @(7, 3574.701286306266, 4293.755365342993),
    @(8, 1547.185698823373, 449.3131405452405)
)

foreach ($row in $driftValues) {
    $r = $row[0]
    $drift.Cells.Item($r, 2).Value = $row[1]
    $drift.Cells.Item($r, 3).Value = $row[2]
}

# ---------------------------------------------------------------------
# 2) Numerical sheet - refresh counts/stats, add "skewness" column (L)
# ---------------------------------------------------------------------
$numerical = $wb.Worksheets.Item("Numerical")

# Copy the header formatting from the existing "Median" header (K1) into
# the new "skewness" header (L1), then set its text.
$numerical.Range("K1").Copy()
$numerical.Range("L1").PasteSpecial(-4122)
$numerical.Range("L1").Value = "skewness"

# row -> B, C, G, I, J, K, L
$numericalValues = @(
    @(2, 108775, 0,     15282, 37126,  17668.61, 16602,  0.1),
    @(3, 108775, 9590,  11121, 195125, 1812.12,  600,    8.4),
    @(4, 108775, 53776, 9179,  115864, 762.04,   21,     9.17),
    @(5, 108775, 7921,  27280, 622800, 10239.19, 3395,   7.15),
    @(6, 108775, 19076, 21008, 470000, 6090.47,  1661.1, 7.61),
    @(7, 108775, 1724,  1876,  648000, 17619.99, 6600,   5.99),
    @(8, 108775, 0,     5323,  15718,  5703.41,  4953,   0.44)
)

foreach ($row in $numericalValues) {
    $r = $row[0]
    $numerical.Cells.Item($r, 2).Value = $row[1]   # Count
    $numerical.Cells.Item($r, 3).Value = $row[2]   # Zeros
    $numerical.Cells.Item($r, 7).Value = $row[3]   # Unique Values
    $numerical.Cells.Item($r, 9).Value = $row[4]   # Maximum
    $numerical.Cells.Item($r, 10).Value = $row[5]  # Mean
    $numerical.Cells.Item($r, 11).Value = $row[6]  # Median
    $numerical.Cells.Item($r, 12).Value = $row[7]  # skewness
}

# ---------------------------------------------------------------------
# 3) Categorical sheet - refresh row count (column B, rows 2-3)
# ---------------------------------------------------------------------
$categorical = $wb.Worksheets.Item("Categorical")
$categorical.Cells.Item(2, 2).Value = 108775
$categorical.Cells.Item(3, 2).Value = 108775

# ---------------------------------------------------------------------
# 4) New "Usage" sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$usage = $wb.Worksheets.Add($null, $lastSheet)
$usage.Name = "Usage"

# Header row, copying the bold/bordered header style used elsewhere
$numerical.Range("A1").Copy()
$usage.Range("A1:F1").PasteSpecial(-4122)

$usage.Range("A1").Value = "schemaname"
$usage.Range("B1").Value = "no_of_times_accessed"
$usage.Range("C1").Value = "table_name"
$usage.Range("D1").Value = "indexrelname"
$usage.Range("E1").Value = "tables_usability"
$usage.Range("F1").Value = "index_usability"

$usage.Range("A2").Value = "adaptiveai"
$usage.Range("B2").Value = 28
$usage.Range("C2").Value = "project_month_dim"
$usage.Range("E2").Value = "Used"
$usage.Range("F2").Value = "Index not used"
